$d = $word.ActiveDocument

# The document ends with a trailing empty paragraph. Insert the two new
# "DICAS DE PREPARO" paragraphs right before that trailing empty paragraph,
# i.e. right after the "OBSERVAÇÃO: ..." paragraph.

$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing.Range.InsertParagraphBefore()
$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara1.Range.Text = "DICAS DE PREPARO: - Para um bolo mais fofo, peneire a farinha de trigo. "

$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing.Range.InsertParagraphBefore()
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara2.Range.Text = "- Você pode adicionar nozes picadas à massa para um toque especial. "
